# Insert a new row at position 340, shifting existing rows 340-355 down to 341-356.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with the new weekly price record.
$ws.Range("A340").Value = 5
$ws.Range("B340").Value = "Macroferia Regional de Talca"
$ws.Range("C340").Value = "Maule"
$ws.Range("D340").Value = 44753
$ws.Range("E340").Value = 7
$ws.Range("F340").Value = 100112023
$ws.Range("G340").Value = "Brócoli"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 5000
$ws.Range("K340").Value = 900
$ws.Range("L340").Value = 900
$ws.Range("M340").Value = 900
$ws.Range("N340").Value = "$/unidad"
$ws.Range("O340").Value = "Región del Maule"
$ws.Range("P340").Value = 900
$ws.Range("Q340").Value = 1
$ws.Range("R340").Value = "Hortaliza"
